# Auto-generated cryptos list update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'39.406.36"
$ws.Range("E2").Value = "  +1.60%  "
# Row 3
$ws.Range("D3").Value = "'2.161.41"
$ws.Range("E3").Value = "  +3.28%  "
# Row 5
$ws.Range("D5").Value = "'228.50"
$ws.Range("E5").Value = "  -0.28%  "
# Row 6
$ws.Range("E6").Value = "  +0.98%  "
# Row 7
$ws.Range("D7").Value = "'64.34"
$ws.Range("E7").Value = "  +5.18%  "
# Row 8
$ws.Range("E8").Value = "  +0.03%  "
# Row 9
$ws.Range("E9").Value = "  +2.68%  "
# Row 10
$ws.Range("E10").Value = "  +1.96%  "
# Row 11
$ws.Range("E11").Value = "  +0.01%  "
# Row 12
$ws.Range("D12").Value = "'15.99"
$ws.Range("E12").Value = "  +4.26%  "
# Row 13
$ws.Range("D13").Value = "'2.481.43"
$ws.Range("E13").Value = "  +3.20%  "
# Row 14
$ws.Range("D14").Value = "'22.33"
$ws.Range("E14").Value = "  +1.08%  "
# Row 15
$ws.Range("D15").Value = "'0.814"
# Row 16
$ws.Range("D16").Value = "'5.56"
$ws.Range("E16").Value = "  +1.22%  "
# Row 17
$ws.Range("D17").Value = "'2.148.50"
$ws.Range("E17").Value = "  +3.05%  "
# Row 18
$ws.Range("D18").Value = "'39.432.73"
$ws.Range("E18").Value = "  +1.88%  "
# Row 19
$ws.Range("B19").Value = "Litecoin"
$ws.Range("C19").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D19").Value = "'71.90"
$ws.Range("E19").Value = "  +0.03%  "
# Row 20
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'6.14"
$ws.Range("E20").Value = "  +0.80%  "
# Row 21
$ws.Range("D21").Value = "'0.0₃0856"
$ws.Range("E21").Value = "  +1.72%  "
# Row 22
$ws.Range("D22").Value = "'231.55"
$ws.Range("E22").Value = "  +1.61%  "
# Row 23
$ws.Range("E23").Value = "  -0.02%  "
# Row 24
$ws.Range("D24").Value = "'2.52"
$ws.Range("E24").Value = "  +5.58%  "
# Row 25
$ws.Range("E25").Value = "  +0.86%  "
# Row 26
$ws.Range("D26").Value = "'172.24"
$ws.Range("E26").Value = "  +0.52%  "
# Row 27
$ws.Range("D27").Value = "'9.54"
$ws.Range("E27").Value = "  -0.01%  "
# Row 28
$ws.Range("E28").Value = "  +0.67%  "
# Row 29
$ws.Range("D29").Value = "'19.94"
$ws.Range("E29").Value = "  +2.91%  "
# Row 30
$ws.Range("E30").Value = "  -1.58%  "
# Row 31
$ws.Range("D31").Value = "'2.69"
$ws.Range("E31").Value = "  +9.14%  "
# Row 32
$ws.Range("E32").Value = "  +1.18%  "
# Row 33
$ws.Range("D33").Value = "'4.63"
$ws.Range("E33").Value = "  +2.35%  "
# Row 34
$ws.Range("D34").Value = "'4.79"
$ws.Range("E34").Value = "  +1.01%  "
# Row 35
$ws.Range("E35").Value = "  +9.92%  "
# Row 36
$ws.Range("D36").Value = "'0.0619"
$ws.Range("E36").Value = "  +1.20%  "
# Row 37
$ws.Range("E37").Value = "  +0.48%  "
# Row 38
$ws.Range("D38").Value = "'3.61"
$ws.Range("E38").Value = "  +0.52%  "
# Row 39
$ws.Range("E39").Value = "  -0.05%  "
# Row 40
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.0232"
$ws.Range("E40").Value = "  +1.36%  "
# Row 41
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'103.88"
$ws.Range("E41").Value = "  +2.97%  "
# Row 42
$ws.Range("D42").Value = "'17.85"
$ws.Range("E42").Value = "  -0.79%  "
# Row 43
$ws.Range("D43").Value = "'1.539.70"
$ws.Range("E43").Value = "  +0.45%  "
# Row 44
$ws.Range("E44").Value = "  +4.51%  "
# Row 45
$ws.Range("D45").Value = "'7.97"
$ws.Range("E45").Value = "  +3.39%  "
# Row 46
$ws.Range("D46").Value = "'4.31"
$ws.Range("E46").Value = "  +4.42%  "
# Row 47
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.0927"
$ws.Range("E47").Value = "  +1.45%  "
# Row 48
$ws.Range("B48").Value = "HuobiToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D48").Value = "'2.82"
$ws.Range("E48").Value = "  +0.54%  "
# Row 49
$ws.Range("E49").Value = "  +5.78%  "
# Row 50
$ws.Range("D50").Value = "'2.365.29"
# Row 51
$ws.Range("E51").Value = "  -0.33%  "

Write-Host "Updated 92 cells"
